$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 21 (value 934), shifting subsequent rows up.
$ws.Rows.Item(21).Delete()
